# Update the date/title line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-11 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-12 Friday", 2)

# Update the division problems in the first (only) table.
# Rows 1, 5, 9, 13, 17 (1-indexed) hold the five problem rows, each with 5 cells.
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("94÷9=", "74÷4=", "58÷7=", "75÷7=", "70÷3=")
    5  = @("51÷4=", "76÷3=", "22÷7=", "41÷6=", "87÷5=")
    9  = @("43÷4=", "32÷3=", "22÷6=", "16÷4=", "43÷7=")
    13 = @("33÷4=", "23÷7=", "89÷4=", "43÷2=", "65÷4=")
    17 = @("77÷5=", "16÷5=", "31÷8=", "31÷2=", "33÷6=")
}

foreach ($rowIndex in $newValues.Keys) {
    $vals = $newValues[$rowIndex]
    for ($col = 1; $col -le $vals.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $vals[$col - 1]
    }
}
